$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Price cells are forced to Text so values like "314.96" or "45.332.66"
# (thousands-grouped, not true decimals) are stored verbatim, matching the
# original inline-string cells, instead of being auto-coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.332.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.373.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.636'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.98%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.26%  '
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.983'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.734.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.369.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.298.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +20.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.05%  '
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0972'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '166.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.18%  '
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.118'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("E39").Value = '  +7.71%  '
$ws.Range("E40").Value = '  -7.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.227'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.75%  '
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.812.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.63%  '
